$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet is protected; unprotect temporarily so the values can be updated.
$ws.Unprotect()

# Update the confidentiality footer text (date 2021-03-31 -> 2021-04-05)
$ws.Range("A58").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-05 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-55
$ws.Range("D2").Value = 0.01641206345319795
$ws.Range("E2").Value = -0.002487167275228797
$ws.Range("D3").Value = 0.05083947537778748
$ws.Range("E3").Value = 0.02079405251502697
$ws.Range("D4").Value = 0.01467766695025905
$ws.Range("E4").Value = 0.0199211045364891
$ws.Range("D5").Value = 0.009558657705323209
$ws.Range("E5").Value = 0.0133261542603309
$ws.Range("D6").Value = 0.01545134051038973
$ws.Range("E6").Value = 0.01330690826727055
$ws.Range("D7").Value = 0.02030092144418078
$ws.Range("E7").Value = 0.01616026213881439
$ws.Range("D8").Value = 0.004431285117395764
$ws.Range("E8").Value = -0.02576219512195144
$ws.Range("D9").Value = 0.006645769674756346
$ws.Range("E9").Value = -0.01550792821048963
$ws.Range("D10").Value = 0.01397012813316964
$ws.Range("E10").Value = 0.0037715517241379
$ws.Range("D11").Value = 0.008812390176841923
$ws.Range("E11").Value = 0.01058548693239914
$ws.Range("D12").Value = 0.01466029693019957
$ws.Range("E12").Value = 0.006582411795681864
$ws.Range("D13").Value = 0.003023991825539894
$ws.Range("E13").Value = 0.04653760238272531
$ws.Range("D14").Value = 0.006054095324804418
$ws.Range("E14").Value = 0.0111046171829341
$ws.Range("D15").Value = 0.01433859772535728
$ws.Range("E15").Value = -0.0005855181835924261
$ws.Range("D16").Value = 0.01058702722625194
$ws.Range("E16").Value = -0.005332239540607131
$ws.Range("D17").Value = 0.02193254532843443
$ws.Range("E17").Value = 0.004289862724392579
$ws.Range("D18").Value = 0.00850185281822303
$ws.Range("E18").Value = 0.003654854600349644
$ws.Range("D19").Value = 0.01676065402242862
$ws.Range("E19").Value = 0.003684824663759745
$ws.Range("D20").Value = 0.01191712042895458
$ws.Range("E20").Value = -0.001889440725545311
$ws.Range("D21").Value = 0.007252884042538921
$ws.Range("E21").Value = 0.007876599934361739
$ws.Range("D22").Value = 0.01338392212286607
$ws.Range("E22").Value = 0.009959623149394359
$ws.Range("D23").Value = 0.01922687520383638
$ws.Range("E23").Value = 0.004155750293612925
$ws.Range("D24").Value = 0.009846067203899942
$ws.Range("E24").Value = -0.01269540502131694
$ws.Range("D25").Value = 0.02115758510015119
$ws.Range("E25").Value = 0.01193770258518456
$ws.Range("D26").Value = 0.01156997736139553
$ws.Range("E26").Value = 0.003911723493694508
$ws.Range("D27").Value = 0.02262059112301265
$ws.Range("E27").Value = 0.03363053339590172
$ws.Range("D28").Value = 0.05539106396744704
$ws.Range("E28").Value = 0.02357723577235782
$ws.Range("D29").Value = 0.02136670084164503
$ws.Range("E29").Value = 0.03387133439418033
$ws.Range("D30").Value = 0.03254498425069931
$ws.Range("E30").Value = 0.01664426346169057
$ws.Range("D31").Value = 0.01639932543848766
$ws.Range("E31").Value = 0.01843384861020403
$ws.Range("D32").Value = 0.01336989744000323
$ws.Range("E32").Value = 0.01499360029255814
$ws.Range("D33").Value = 0.0214021807159517
$ws.Range("E33").Value = 0.01450059517368252
$ws.Range("D34").Value = 0.0411048014691949
$ws.Range("E34").Value = 0.04187286949825797
$ws.Range("D35").Value = 0.01110523282469304
$ws.Range("E35").Value = 0.001390337156760513
$ws.Range("D36").Value = 0.009774045954060737
$ws.Range("E36").Value = -0.02193144120899371
$ws.Range("D37").Value = 0.01178025110422665
$ws.Range("E37").Value = 0.02908048330944091
$ws.Range("D38").Value = 0.007224963343628503
$ws.Range("E38").Value = -0.003339121143315027
$ws.Range("D39").Value = 0.01181344714256254
$ws.Range("E39").Value = 0.007850088630032781
$ws.Range("D40").Value = 0.01772346780098473
$ws.Range("E40").Value = 0.008272412466333146
$ws.Range("D41").Value = 0.01705034735697614
$ws.Range("E41").Value = 0.02057495160943423
$ws.Range("D42").Value = 0.0337325789555066
$ws.Range("E42").Value = 0.01129991989930201
$ws.Range("D43").Value = 0.0112723710177098
$ws.Range("E43").Value = 0.007048362611147274
$ws.Range("D44").Value = 0.02162458130612065
$ws.Range("E44").Value = 0.01005256847735847
$ws.Range("D45").Value = 0.01377876841218106
$ws.Range("E45").Value = 0.02632629477886317
$ws.Range("D46").Value = 0.008208492479440733
$ws.Range("E46").Value = -0.006947873315934272
$ws.Range("D47").Value = 0.01347566156214316
$ws.Range("E47").Value = 0.004167740826673594
$ws.Range("D48").Value = 0.009665547662096628
$ws.Range("E48").Value = -0.0011681193278823
$ws.Range("D49").Value = 0.01479983609134405
$ws.Range("E49").Value = 0.02241696334259213
$ws.Range("D50").Value = 0.008333749624091859
$ws.Range("E50").Value = 0.0003280839895010157
$ws.Range("D51").Value = 0.01112935785255343
$ws.Range("E51").Value = 0.01746293245469532
$ws.Range("D52").Value = 0.008547915538084462
$ws.Range("E52").Value = 0.01452934845599807
$ws.Range("D53").Value = 0.1396835253111912
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 0.04376312220577896
$ws.Range("E54").Value = 0.01399475196801214
$ws.Range("E55").Value = 0.01161591141444185

# Re-apply the original sheet protection so the workbook state matches.
$ws.Protect("D382")
